# Atualizacao de bases das ligas - Poland Ekstraklasa
# Rearranges match-result rows 298-306 (id 296-304) to the updated source order;
# columns A (id) and D (Date) are unchanged, only B,E:AD are rewritten per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 298
$ws.Cells.Item(298, 2).Value = 7083188
$ws.Cells.Item(298, 5).Value = "Legia Warsaw"
$ws.Cells.Item(298, 6).Value = "Zaglebie Lubin"
$ws.Cells.Item(298, 7).Value = 2
$ws.Cells.Item(298, 8).Value = 1
$ws.Cells.Item(298, 9).Value = 2
$ws.Cells.Item(298, 10).Value = 0
$ws.Cells.Item(298, 11).Value = "H"
$ws.Cells.Item(298, 12).Value = 1.5
$ws.Cells.Item(298, 13).Value = 4
$ws.Cells.Item(298, 14).Value = 5.5
$ws.Cells.Item(298, 15).Value = 1.6
$ws.Cells.Item(298, 16).Value = 4.1
$ws.Cells.Item(298, 17).Value = 4.333
$ws.Cells.Item(298, 18).Value = -0.75
$ws.Cells.Item(298, 19).Value = 1.825
$ws.Cells.Item(298, 20).Value = 2.025
$ws.Cells.Item(298, 21).Value = 3
$ws.Cells.Item(298, 22).Value = 1.875
$ws.Cells.Item(298, 23).Value = 1.975
$ws.Cells.Item(298, 24).Value = 0.6000000000000001
$ws.Cells.Item(298, 25).Value = -1
$ws.Cells.Item(298, 26).Value = -1
$ws.Cells.Item(298, 27).Value = 0.4125
$ws.Cells.Item(298, 28).Value = -0.5
$ws.Cells.Item(298, 29).Value = 0
$ws.Cells.Item(298, 30).Value = 0

# Row 299
$ws.Cells.Item(299, 2).Value = 7093820
$ws.Cells.Item(299, 5).Value = "Ruch Chorzow"
$ws.Cells.Item(299, 6).Value = "Cracovia Krakow"
$ws.Cells.Item(299, 7).Value = 2
$ws.Cells.Item(299, 8).Value = 0
$ws.Cells.Item(299, 9).Value = 1
$ws.Cells.Item(299, 10).Value = 0
$ws.Cells.Item(299, 11).Value = "H"
$ws.Cells.Item(299, 12).Value = 2.5
$ws.Cells.Item(299, 13).Value = 3.4
$ws.Cells.Item(299, 14).Value = 2.5
$ws.Cells.Item(299, 15).Value = 2.6
$ws.Cells.Item(299, 16).Value = 3.6
$ws.Cells.Item(299, 17).Value = 2.3
$ws.Cells.Item(299, 18).Value = 0
$ws.Cells.Item(299, 19).Value = 2.025
$ws.Cells.Item(299, 20).Value = 1.825
$ws.Cells.Item(299, 21).Value = 3
$ws.Cells.Item(299, 22).Value = 2.025
$ws.Cells.Item(299, 23).Value = 1.825
$ws.Cells.Item(299, 24).Value = 1.6
$ws.Cells.Item(299, 25).Value = -1
$ws.Cells.Item(299, 26).Value = -1
$ws.Cells.Item(299, 27).Value = 1.025
$ws.Cells.Item(299, 28).Value = -1
$ws.Cells.Item(299, 29).Value = -1
$ws.Cells.Item(299, 30).Value = 0.825

# Row 300
$ws.Cells.Item(300, 2).Value = 7083187
$ws.Cells.Item(300, 5).Value = "Lech Poznan"
$ws.Cells.Item(300, 6).Value = "Korona Kielce"
$ws.Cells.Item(300, 7).Value = 1
$ws.Cells.Item(300, 8).Value = 2
$ws.Cells.Item(300, 9).Value = 1
$ws.Cells.Item(300, 10).Value = 0
$ws.Cells.Item(300, 11).Value = "A"
$ws.Cells.Item(300, 12).Value = 1.8
$ws.Cells.Item(300, 13).Value = 3.8
$ws.Cells.Item(300, 14).Value = 3.6
$ws.Cells.Item(300, 15).Value = 2.1
$ws.Cells.Item(300, 16).Value = 3.7
$ws.Cells.Item(300, 17).Value = 2.9
$ws.Cells.Item(300, 18).Value = -0.25
$ws.Cells.Item(300, 19).Value = 1.9
$ws.Cells.Item(300, 20).Value = 1.95
$ws.Cells.Item(300, 21).Value = 2.75
$ws.Cells.Item(300, 22).Value = 1.925
$ws.Cells.Item(300, 23).Value = 1.925
$ws.Cells.Item(300, 24).Value = -1
$ws.Cells.Item(300, 25).Value = -1
$ws.Cells.Item(300, 26).Value = 1.9
$ws.Cells.Item(300, 27).Value = -1
$ws.Cells.Item(300, 28).Value = 0.95
$ws.Cells.Item(300, 29).Value = 0.4625
$ws.Cells.Item(300, 30).Value = -0.5

# Row 301
$ws.Cells.Item(301, 2).Value = 7041338
$ws.Cells.Item(301, 5).Value = "Jagiellonia Bialystok"
$ws.Cells.Item(301, 6).Value = "Warta Poznan"
$ws.Cells.Item(301, 7).Value = 3
$ws.Cells.Item(301, 8).Value = 0
$ws.Cells.Item(301, 9).Value = 3
$ws.Cells.Item(301, 10).Value = 0
$ws.Cells.Item(301, 11).Value = "H"
$ws.Cells.Item(301, 12).Value = 1.444
$ws.Cells.Item(301, 13).Value = 4.75
$ws.Cells.Item(301, 14).Value = 5.25
$ws.Cells.Item(301, 15).Value = 1.4
$ws.Cells.Item(301, 16).Value = 4.75
$ws.Cells.Item(301, 17).Value = 5.75
$ws.Cells.Item(301, 18).Value = -1.25
$ws.Cells.Item(301, 19).Value = 1.9
$ws.Cells.Item(301, 20).Value = 1.95
$ws.Cells.Item(301, 21).Value = 3
$ws.Cells.Item(301, 22).Value = 1.925
$ws.Cells.Item(301, 23).Value = 1.925
$ws.Cells.Item(301, 24).Value = 0.3999999999999999
$ws.Cells.Item(301, 25).Value = -1
$ws.Cells.Item(301, 26).Value = -1
$ws.Cells.Item(301, 27).Value = 0.8999999999999999
$ws.Cells.Item(301, 28).Value = -1
$ws.Cells.Item(301, 29).Value = 0
$ws.Cells.Item(301, 30).Value = 0

# Row 302
$ws.Cells.Item(302, 2).Value = 7090293
$ws.Cells.Item(302, 5).Value = "Radomiak Radom"
$ws.Cells.Item(302, 6).Value = "Widzew Lodz"
$ws.Cells.Item(302, 7).Value = 1
$ws.Cells.Item(302, 8).Value = 3
$ws.Cells.Item(302, 9).Value = 1
$ws.Cells.Item(302, 10).Value = 0
$ws.Cells.Item(302, 11).Value = "A"
$ws.Cells.Item(302, 12).Value = 2.2
$ws.Cells.Item(302, 13).Value = 3.1
$ws.Cells.Item(302, 14).Value = 3.1
$ws.Cells.Item(302, 15).Value = 2.15
$ws.Cells.Item(302, 16).Value = 3.2
$ws.Cells.Item(302, 17).Value = 3.1
$ws.Cells.Item(302, 18).Value = -0.25
$ws.Cells.Item(302, 19).Value = 1.925
$ws.Cells.Item(302, 20).Value = 1.925
$ws.Cells.Item(302, 21).Value = 2.75
$ws.Cells.Item(302, 22).Value = 1.9
$ws.Cells.Item(302, 23).Value = 1.95
$ws.Cells.Item(302, 24).Value = -1
$ws.Cells.Item(302, 25).Value = -1
$ws.Cells.Item(302, 26).Value = 2.1
$ws.Cells.Item(302, 27).Value = -1
$ws.Cells.Item(302, 28).Value = 0.925
$ws.Cells.Item(302, 29).Value = 0.8999999999999999
$ws.Cells.Item(302, 30).Value = -1

# Row 303
$ws.Cells.Item(303, 2).Value = 7083189
$ws.Cells.Item(303, 5).Value = "Pogon Szczecin"
$ws.Cells.Item(303, 6).Value = "Gornik Zabrze"
$ws.Cells.Item(303, 7).Value = 1
$ws.Cells.Item(303, 8).Value = 0
$ws.Cells.Item(303, 9).Value = 0
$ws.Cells.Item(303, 10).Value = 0
$ws.Cells.Item(303, 11).Value = "H"
$ws.Cells.Item(303, 12).Value = 1.727
$ws.Cells.Item(303, 13).Value = 4
$ws.Cells.Item(303, 14).Value = 3.75
$ws.Cells.Item(303, 15).Value = 1.55
$ws.Cells.Item(303, 16).Value = 4.333
$ws.Cells.Item(303, 17).Value = 4.5
$ws.Cells.Item(303, 18).Value = -1
$ws.Cells.Item(303, 19).Value = 1.925
$ws.Cells.Item(303, 20).Value = 1.925
$ws.Cells.Item(303, 21).Value = 3.5
$ws.Cells.Item(303, 22).Value = 2.025
$ws.Cells.Item(303, 23).Value = 1.825
$ws.Cells.Item(303, 24).Value = 0.55
$ws.Cells.Item(303, 25).Value = -1
$ws.Cells.Item(303, 26).Value = -1
$ws.Cells.Item(303, 27).Value = 0
$ws.Cells.Item(303, 28).Value = 0
$ws.Cells.Item(303, 29).Value = -1
$ws.Cells.Item(303, 30).Value = 0.825

# Row 304
$ws.Cells.Item(304, 2).Value = 7088350
$ws.Cells.Item(304, 5).Value = "Puszcza Niepolomice"
$ws.Cells.Item(304, 6).Value = "Piast Gliwice"
$ws.Cells.Item(304, 7).Value = 1
$ws.Cells.Item(304, 8).Value = 0
$ws.Cells.Item(304, 9).Value = 0
$ws.Cells.Item(304, 10).Value = 0
$ws.Cells.Item(304, 11).Value = "H"
$ws.Cells.Item(304, 12).Value = 3
$ws.Cells.Item(304, 13).Value = 3.1
$ws.Cells.Item(304, 14).Value = 2.3
$ws.Cells.Item(304, 15).Value = 2.7
$ws.Cells.Item(304, 16).Value = 3
$ws.Cells.Item(304, 17).Value = 2.625
$ws.Cells.Item(304, 18).Value = 0
$ws.Cells.Item(304, 19).Value = 1.975
$ws.Cells.Item(304, 20).Value = 1.875
$ws.Cells.Item(304, 21).Value = 2.25
$ws.Cells.Item(304, 22).Value = 2.025
$ws.Cells.Item(304, 23).Value = 1.825
$ws.Cells.Item(304, 24).Value = 1.7
$ws.Cells.Item(304, 25).Value = -1
$ws.Cells.Item(304, 26).Value = -1
$ws.Cells.Item(304, 27).Value = 0.9750000000000001
$ws.Cells.Item(304, 28).Value = -1
$ws.Cells.Item(304, 29).Value = -1
$ws.Cells.Item(304, 30).Value = 0.825

# Row 305
$ws.Cells.Item(305, 2).Value = 7074364
$ws.Cells.Item(305, 5).Value = "Rakow Czestochowa"
$ws.Cells.Item(305, 6).Value = "Slask Wroclaw"
$ws.Cells.Item(305, 7).Value = 1
$ws.Cells.Item(305, 8).Value = 2
$ws.Cells.Item(305, 9).Value = 1
$ws.Cells.Item(305, 10).Value = 0
$ws.Cells.Item(305, 11).Value = "A"
$ws.Cells.Item(305, 12).Value = 2.5
$ws.Cells.Item(305, 13).Value = 3.6
$ws.Cells.Item(305, 14).Value = 2.4
$ws.Cells.Item(305, 15).Value = 2.15
$ws.Cells.Item(305, 16).Value = 3.6
$ws.Cells.Item(305, 17).Value = 2.875
$ws.Cells.Item(305, 18).Value = -0.25
$ws.Cells.Item(305, 19).Value = 1.95
$ws.Cells.Item(305, 20).Value = 1.9
$ws.Cells.Item(305, 21).Value = 2.5
$ws.Cells.Item(305, 22).Value = 1.875
$ws.Cells.Item(305, 23).Value = 1.975
$ws.Cells.Item(305, 24).Value = -1
$ws.Cells.Item(305, 25).Value = -1
$ws.Cells.Item(305, 26).Value = 1.875
$ws.Cells.Item(305, 27).Value = -1
$ws.Cells.Item(305, 28).Value = 0.8999999999999999
$ws.Cells.Item(305, 29).Value = 0.875
$ws.Cells.Item(305, 30).Value = -1

# Row 306
$ws.Cells.Item(306, 2).Value = 7093821
$ws.Cells.Item(306, 5).Value = "LKS Lodz"
$ws.Cells.Item(306, 6).Value = "Stal Mielec"
$ws.Cells.Item(306, 7).Value = 3
$ws.Cells.Item(306, 8).Value = 2
$ws.Cells.Item(306, 9).Value = 3
$ws.Cells.Item(306, 10).Value = 0
$ws.Cells.Item(306, 11).Value = "H"
$ws.Cells.Item(306, 12).Value = 2.5
$ws.Cells.Item(306, 13).Value = 3.4
$ws.Cells.Item(306, 14).Value = 2.5
$ws.Cells.Item(306, 15).Value = 2.2
$ws.Cells.Item(306, 16).Value = 3.5
$ws.Cells.Item(306, 17).Value = 2.8
$ws.Cells.Item(306, 18).Value = -0.25
$ws.Cells.Item(306, 19).Value = 2.025
$ws.Cells.Item(306, 20).Value = 1.825
$ws.Cells.Item(306, 21).Value = 3
$ws.Cells.Item(306, 22).Value = 2
$ws.Cells.Item(306, 23).Value = 1.85
$ws.Cells.Item(306, 24).Value = 1.2
$ws.Cells.Item(306, 25).Value = -1
$ws.Cells.Item(306, 26).Value = -1
$ws.Cells.Item(306, 27).Value = 1.025
$ws.Cells.Item(306, 28).Value = -1
$ws.Cells.Item(306, 29).Value = 1
$ws.Cells.Item(306, 30).Value = -1

